$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows appended to the BIIB noun-trade data sheet.
$rows = @(
    @{ Row=10; A=9606.7199999999993; B=9524.81;  C=305.24;              D=307.86; E=$false; F=0.86;  G=42613.765462962961; H=$true  },
    @{ Row=11; A=9542.35;             B=9606.7199999999993; C=307.68;   D=305.63; E=$false; F=-0.67; G=42614.672638888886; H=$false },
    @{ Row=12; A=9482.23;             B=9542.35;             C=307.95999999999998; D=306.02; E=$false; F=-0.63; G=42615.750034722223; H=$false }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item(3, 7).Copy()
    $ws.Cells.Item($r.Row, 7).PasteSpecial(-4122)
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}
